$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PC")

# Mark Giuliano Antoniol (row 7) and Gordon Fraser (row 10) as accepted ("Y")
$ws.Range("C7").Value = "Y"
$ws.Range("C10").Value = "Y"

# Update the selected cell to reflect where editing left off
$ws.Range("C13").Select()
